# ----------------------------------------------------------------------------
# Adds "Table S1" (full list of macroeconomic predictors) after "Table 1",
# and "Table Sx - CCI bands" (ESA CCI land-cover band legend) after "Table 2".
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert "Table S1" right after "Table 1"
# ------------------------------------------------------------------
$table1 = $wb.Worksheets.Item("Table 1")
$wsS1 = $wb.Worksheets.Add($null, $table1)
$wsS1.Name = "Table S1"

# Tab-separated rows: Predictor variable | Units | Resolution | Source | Details
$s1data = @"
Predictor variable`tUnits`tResolution`tSource`tDetails
Economy
GDP per capita `tBillions USD`tNational`tWorld Bank`tConstant 2010 rates
GPD growth`t%`tNational`tWorld Bank`tConstant 2010 rates
GNI per capita`tUSD `tNational`tWorld Bank`tGross National Income per capita. Calculated as gross national income divded by the mid-year population at current USD rates
Foreign Direct Investment`tMillions USD`tNational`tUNCTAD`tInward and outward flows and stock
Agricultural sector value of GDP`t%`tNational`tCNIS`tProportion of national GDP
Industrial sector value of GDP`t%`tNational`tCNIS`tProportion of national GDP
Development flows to agriculture`tMillions USD`tNational`tFAO`tDonor development investment flows, other official flows, and private donor flows at constant 2016 prices to all agriculture and forestry sub-sectors
Development flows to environment`tMillions USD`tNational`tFAO`tDonor development investment flows, other official flows, and private donor flows at constant 2016 prices to general environment protection
Population density`tpax/km2`tNational`tFAO
Commodity prices
Agricultural Raw Materials `tIndex`tGlobal`tIMF`tPrice index for global agricultural raw materials including timber, cotton, wool, rubber, and hides
Crop Production `tIndex`tNational`tFAO`tRelative level of the aggregate volume of agricultural production for each year in comparison with the base period 2004-2006
Non-food agricultural production `tIndex`tNational`tFAO`tRelative level of the aggregate volume of non-food agricultural production for each year in comparison with the base period 2004-2006
Forestry production`tm3`tNational`tFAO`tTotal production values for industrial roundwood, non-coniferous tropical wood, other industrial roundwood, sawlogs and veneer logs (coniferous and non-coniferous), and sawnwood (coniferous and non-coniferous
Price of rice`tUSD/ton`tGlobal`tWorld Bank`tMedian annual global market price of rice
Price of corn`tUSD/ton`tGlobal`tWorld Bank`tAnnual global market price of corn
Price of rubber`tUSD/ton`tRegional`tRASCE`tMonthly regional market value of rubber on the Singapore Exchange
Price of sugar`tUSD/ton`tGlobal`tWorld Bank`tAnnual global market price of sugar
Producer prices
Producer price of Rice`tUSD/ton`tNational`tFAO`tFarmgate prices for Cambodian producers
Producer price of rubber`tUSD/ton`tNational`tFAO`tFarmgate prices for Cambodian producers
Producer price of cassava`tUSD/ton`tNational`tFAO`tFarmgate prices for Cambodian producers
Producer price of corn`tUSD/ton`tNational`tFAO`tFarmgate prices for Cambodian producers
Producer price of sugar`tUSD/ton`tNational`tFAO`tFarmgate prices for Cambodian producers
Control
Forest remaining`tkm2`tNational`tESACCI`tTotal forested area
"@

$s1lines = $s1data -split "`r?`n"
$s1rows = $s1lines.Count
$arrS1 = New-Object 'object[,]' $s1rows,5
for ($r = 0; $r -lt $s1rows; $r++) {
    $parts = $s1lines[$r] -split "`t"
    for ($c = 0; $c -lt $parts.Length; $c++) {
        if ($parts[$c].Length -gt 0) {
            $arrS1[$r,$c] = $parts[$c]
        }
    }
}
$wsS1.Range("A1:E$s1rows").Value = $arrS1

# Header row: bold
$wsS1.Range("A1:E1").Font.Bold = $true
# Section header rows (Economy / Commodity prices / Producer prices / Control): italic
$wsS1.Range("A2").Font.Italic = $true
$wsS1.Range("A12").Font.Italic = $true
$wsS1.Range("A21").Font.Italic = $true
$wsS1.Range("A27").Font.Italic = $true
# Sub-group label cells that repeat the (normal) font explicitly
$wsS1.Range("A13").Font.Italic = $false
$wsS1.Range("A22").Font.Italic = $false

# Wrap text for the "Details" column entries that have explanatory text
$wsS1.Range("E5").WrapText = $true
$wsS1.Range("E9").WrapText = $true
$wsS1.Range("E10").WrapText = $true
$wsS1.Range("E13").WrapText = $true
$wsS1.Range("E14").WrapText = $true
$wsS1.Range("E15").WrapText = $true
$wsS1.Range("E16").WrapText = $true
$wsS1.Range("E17").WrapText = $true
$wsS1.Range("E18").WrapText = $true
$wsS1.Range("E19").WrapText = $true
$wsS1.Range("E20").WrapText = $true
$wsS1.Range("E22").WrapText = $true
$wsS1.Range("E23").WrapText = $true
$wsS1.Range("E24").WrapText = $true
$wsS1.Range("E25").WrapText = $true
$wsS1.Range("E26").WrapText = $true
$wsS1.Range("E28").WrapText = $true

# Row heights (explicit, matching the authored workbook)
$wsS1.Rows.Item(5).RowHeight = 28.8
$wsS1.Rows.Item(9).RowHeight = 33.6
$wsS1.Rows.Item(10).RowHeight = 28.8
$wsS1.Rows.Item(13).RowHeight = 28.8
$wsS1.Rows.Item(14).RowHeight = 28.8
$wsS1.Rows.Item(15).RowHeight = 28.8
$wsS1.Rows.Item(16).RowHeight = 43.2

# Column widths
$wsS1.Columns.Item(1).ColumnWidth = 31.6
$wsS1.Columns.Item(2).ColumnWidth = 10.94
$wsS1.Columns.Item(3).ColumnWidth = 9.39
$wsS1.Columns.Item(4).ColumnWidth = 9.94
$wsS1.Columns.Item(5).ColumnWidth = 73.17

$wsS1.Range("E28").Select()

# ------------------------------------------------------------------
# 2. Insert "Table Sx - CCI bands" right after "Table 2"
# ------------------------------------------------------------------
$table2 = $wb.Worksheets.Item("Table 2")
$wsSx = $wb.Worksheets.Add($null, $table2)
$wsSx.Name = "Table Sx - CCI bands"

$sxdata = @"
Value`tLabel
0`tNo data
10`tCropland, rainfed
11`tHerbaceous cover
12`tTree or shrub cover
20`tCropland, irrigated or post-flooding
30`tMosaic cropland (>50%) / natural vegetation (tree, shrub, herbaceous cover) (<50%)
40`tMosaic natural vegetation (tree, shrub, herbaceous cover) (>50%) / cropland (<50%)
50`tTree cover, broadleaved, evergreen, cosed to open (>15%)
60`tTree cover, broadleaved, deciduous, closed to open (>15%)
61`tTree cover, broadleaves, decisuous, closed (>40%)
62`tTree cover, broadleaves, deciduous, open (15 - 40%)
70`tTree cover, needleleaved, evergreen, closed to open (>15%)
71`tTree cover, needleleaved, evergreen, closed (>40%)
72`tTree cover, needleleaved, evergreen, open (15 - 40%)
80`tTree cover, needleleaved, deciduous, closed to open (>15%)
81`tTree cover, needleleaved, deciduous, closed (>40%)
82`tTree cover, needleleaved, deciduous, open (15 - 40%)
90`tTree cover, mixed leaf type (broadleaved and needleleaved)
100`tMosaic tree and shrub (>50%) / herbaceous cover (<50%)
110`tMosaic herbaceous cover (>50%) / tree and shrub (<50%)
120`tShrubland
121`tEvergreen shrubland
122`tDeciduous shrubland
130`tGrassland
140`tLichens and mosses
150`tSparse vegetation (tree, shrub, herbaceous cover) (<15%)
152`tSparse shrub (<15%)
153`tSparse herbaceous cover (<15%)
160`tTree cover, flooded, fresh or brakish water
"@

$sxlines = $sxdata -split "`r?`n"
$sxrows = $sxlines.Count
$arrSx = New-Object 'object[,]' $sxrows,2
for ($r = 0; $r -lt $sxrows; $r++) {
    $parts = $sxlines[$r] -split "`t"
    $arrSx[$r,1] = $parts[1]
    if ($r -eq 0) {
        $arrSx[$r,0] = $parts[0]
    } else {
        $arrSx[$r,0] = [double]$parts[0]
    }
}
$wsSx.Range("A1:B$sxrows").Value = $arrSx

# Green fill for the CCI "tree cover" band rows (A9:B20 -> value rows 50-100)
$wsSx.Range("A9:B20").Interior.Color = 5296274

$wsSx.Columns.Item(2).ColumnWidth = 73.27

$wsSx.Range("B2").Select()

# ------------------------------------------------------------------
# 3. Restore "Table 1" as the active/selected sheet & cell, matching
#    the authored window state.
# ------------------------------------------------------------------
$table1.Activate()
$table1.Range("A8").Select()

$win = $wb.Windows.Item(1)
$win.Left = 28680
$win.Top = -120
$win.Width = 19440
$win.Height = 15000
